$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, columns B..K hold datetime serial numbers (date + fractional time
# of day). Truncate each one down to the whole-number date serial (drop the
# time-of-day portion), matching the already date-only values used
# elsewhere in the sheet (e.g. row 3).
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
foreach ($col in $cols) {
    $addr = $col + "2"
    $cell = $ws.Range($addr)
    $value = $cell.Value2
    $cell.Value = [math]::Floor($value)
}
